$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name
$ws.Name = "Through 2022-06-20"

# Update the "June" row label (month through date)
$ws.Range("A7").Value = "June (through 06-20)"

# Update June row (row 7) values
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = 23
$ws.Range("D7").Value = 48
$ws.Range("E7").Value = 37
$ws.Range("F7").Value = 34
$ws.Range("G7").Value = 80
$ws.Range("H7").Value = 78
$ws.Range("I7").Value = 99

# Update Total row (row 8) values
$ws.Range("B8").Value = 119
$ws.Range("C8").Value = 232
$ws.Range("D8").Value = 364
$ws.Range("E8").Value = 332
$ws.Range("F8").Value = 238
$ws.Range("G8").Value = 438
$ws.Range("H8").Value = 709
$ws.Range("I8").Value = 762
